# Weekly update: insert a new week's price block for "Pimiento" (Agrícola del
# Norte S.A. de Arica) at the top of its sub-table (rows 834-860), pushing the
# existing rows down by 3 to make room, and append the matching new
# "Zafiro verde" trio that used to be missing at the bottom of the block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that vary row-to-row within this block.
$varyingCols = 4, 8, 9, 10, 11, 12, 13, 16   # D,H,I,J,K,L,M,P

# Capture the date cell's number format so new cells we create keep the
# same "$/caja" style date formatting as the rest of the column.
$dateFormat = $ws.Cells.Item(834, 4).NumberFormat

# Shift the existing block (rows 834-860) down by 3 rows -> (837-863).
# Walk bottom-up so we never overwrite a source row before it's read.
for ($r = 860; $r -ge 834; $r--) {
    $destRow = $r + 3
    foreach ($c in $varyingCols) {
        $val = $ws.Cells.Item($r, $c).Value2
        $destCell = $ws.Cells.Item($destRow, $c)
        $destCell.Value = $val
        if ($c -eq 4) {
            $destCell.NumberFormat = $dateFormat
        }
    }
}

# Values for the new week's block (rows 834-836) and the newly-appended
# "Zafiro verde" trio (rows 861-863), keyed by destination row.
# Columns: D, H, I, J, K, L, M, P
$newData = @{
    834 = @(44939, "Zafiro rojo", "Primera", 400, 7000, 8000, 7375, 492)
    835 = @(44939, "Zafiro rojo", "Segunda", 300, 5000, 6000, 5333, 356)
    836 = @(44939, "Zafiro rojo", "Tercera", 200, 3000, 4000, 3250, 217)
    861 = @(44306, "Zafiro verde", "Primera", 120, 8000, 9000, 8500, 567)
    862 = @(44306, "Zafiro verde", "Segunda", 160, 6000, 7000, 6500, 433)
    863 = @(44306, "Zafiro verde", "Tercera", 160, 5000, 6000, 5500, 367)
}

# Columns that are constant for every row in this block.
$constCells = @{
    1  = 1                                          # A Mercado ID
    2  = "Agrícola del Norte S.A. de Arica"          # B Mercado
    3  = "Arica y Parinacota"                        # C Región
    5  = 15                                          # E Codreg
    6  = 100112002                                   # F Categoría ID
    7  = "Pimiento"                                  # G Categoría
    14 = "`$/caja 15 kilos"                          # N Unidad de comercialización
    15 = "Región de Arica y Parinacota"               # O Origen
    17 = 15                                          # Q Kg o Unidades
    18 = "Hortaliza"                                 # R Clasificación
}

foreach ($destRow in 834, 835, 836, 861, 862, 863) {
    $vals = $newData[$destRow]

    foreach ($colIdx in $constCells.Keys) {
        $ws.Cells.Item($destRow, $colIdx).Value = $constCells[$colIdx]
    }

    $dCell = $ws.Cells.Item($destRow, 4)
    $dCell.Value = $vals[0]
    $dCell.NumberFormat = $dateFormat

    $ws.Cells.Item($destRow, 8).Value = $vals[1]
    $ws.Cells.Item($destRow, 9).Value = $vals[2]
    $ws.Cells.Item($destRow, 10).Value = $vals[3]
    $ws.Cells.Item($destRow, 11).Value = $vals[4]
    $ws.Cells.Item($destRow, 12).Value = $vals[5]
    $ws.Cells.Item($destRow, 13).Value = $vals[6]
    $ws.Cells.Item($destRow, 16).Value = $vals[7]
}

Write-Output "edit complete"
